$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(12318, "Vicente Nascimento", "Recursos Humanos", "Outros", 3, 45080, 4890.74),
    @(64167, "Apollo Teixeira", "Financeiro", "Doenca", 2, 45091, 4865.04),
    @(19393, "Maria Clara Rocha", "Recursos Humanos", "Viagem de negocios", 7, 45088, 5399.82),
    @(44503, "Sr. Danilo Vasconcelos", "Juridico", "Consulta medica", 8, 45089, 4728.77),
    @(2409, "Joaquim Nascimento", "P&D", "Viagem de negocios", 3, 45103, 6512.64),
    @(55286, "Luísa Vieira", "Juridico", "Problemas pessoais", 5, 45093, 5960.63),
    @(99111, "Yago Cirino", "Engenharia", "Outros", 5, 45083, 2865.01),
    @(26570, "Esther Pinto", "TI", "Doenca", 6, 45104, 8679.780000000001),
    @(74015, "Nicolas Camargo", "Operacoes", "Outros", 5, 45097, 4430.13),
    @(69816, "Kamilly Cavalcante", "Operacoes", "Doenca", 2, 45100, 4112.75)
)

$rowIndex = 2
foreach ($record in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $record[0]
    $ws.Cells.Item($rowIndex, 2).Value = $record[1]
    $ws.Cells.Item($rowIndex, 3).Value = $record[2]
    $ws.Cells.Item($rowIndex, 4).Value = $record[3]
    $ws.Cells.Item($rowIndex, 5).Value = $record[4]
    $ws.Cells.Item($rowIndex, 6).Value = $record[5]
    $ws.Cells.Item($rowIndex, 7).Value = $record[6]
    $rowIndex++
}

$wb.Save()
